# The November 2022 MGU dataset had its timestamp column (column A)
# shifted back by 3 hours (10800 seconds) for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 365
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $cell.Value2 = $old - 10800
}
